# Generate Report for Handback
#
# This script updates the localization-status workbook so the "Overview",
# "zh-cn" and "de-de" sheets reflect that the two files have been handed
# back and are in sync with en-US:
#   - Status text "Ready for handoff" -> "Handed back: in sync with en-US"
#     (Overview!E2:F3, zh-cn!C2:C3, de-de!C2:C3) + widen those columns so
#     the longer text fits.
#   - Fill in the "Latest Target File" / "Latest Handback File" /
#     "Latest Handback DateTime" columns (I/J/K) on the zh-cn and de-de
#     detail sheets, including a hyperlink on the new Target File cell
#     that mirrors the Source File Name hyperlink in column A, and widen
#     the I/J columns to fit the long file names.

$wb = $excel.ActiveWorkbook

$oldStatus = "Ready for handoff"
$newStatus = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------
# Overview sheet: status text for each language column + column widths.
# ---------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = $newStatus
$overview.Range("F2").Value = $newStatus
$overview.Range("E3").Value = $newStatus
$overview.Range("F3").Value = $newStatus
$overview.Columns.Item(5).ColumnWidth = 29.1
$overview.Columns.Item(6).ColumnWidth = 29.1

# ---------------------------------------------------------------------
# Helper data per detail sheet: zh-cn / de-de.
# ---------------------------------------------------------------------
$sheetNames = @("zh-cn", "de-de")
$handbackDateTimes = @("2016-09-07 15:44:49", "2016-09-07 15:45:23")

foreach ($idx in 0..1) {
    $sheetName = $sheetNames[$idx]
    $ws = $wb.Worksheets.Item($sheetName)

    # Status column (C) for both data rows.
    $ws.Range("C2").Value = $newStatus
    $ws.Range("C3").Value = $newStatus

    # Look up the existing hyperlinks on column A (Source File Name) so the
    # new column-I (Latest Target File) hyperlinks can mirror them exactly.
    $urlRow2 = ""
    $urlRow3 = ""
    $displayRow2 = ""
    $displayRow3 = ""
    foreach ($h in $ws.Hyperlinks) {
        $addr = $h.Range().Address()
        if ($addr -eq '$A$2') {
            $urlRow2 = $h.Address()
            $displayRow2 = $h.TextToDisplay()
        }
        if ($addr -eq '$A$3') {
            $urlRow3 = $h.Address()
            $displayRow3 = $h.TextToDisplay()
        }
    }

    # Latest Target File (I) now mirrors the Source File Name, with the
    # same hyperlink + "HyperLink" cell style used in column A.
    $ws.Range("I2").Value = $displayRow2
    $ws.Range("I2").Style = "HyperLink"
    $ws.Hyperlinks.Add($ws.Range("I2"), $urlRow2, "", "", $displayRow2)

    $ws.Range("I3").Value = $displayRow3
    $ws.Range("I3").Style = "HyperLink"
    $ws.Hyperlinks.Add($ws.Range("I3"), $urlRow3, "", "", $displayRow3)

    # Latest Handback File (J) = Latest Handoff File (G), now that the
    # handback xlf is in sync with the handoff xlf.
    $ws.Range("J2").Value = $ws.Range("G2").Value()
    $ws.Range("J3").Value = $ws.Range("G3").Value()

    # Latest Handback DateTime (K).
    $ws.Range("K2").Value = $handbackDateTimes[$idx]
    $ws.Range("K3").Value = $handbackDateTimes[$idx]

    # Widen the Status / Target File / Handback File columns to fit the
    # longer text now stored in them.
    $ws.Columns.Item(3).ColumnWidth = 29.1
    $ws.Columns.Item(9).ColumnWidth = 40
    $ws.Columns.Item(10).ColumnWidth = 40
}
